$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 31 with the incoming "Retour" mail log entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A31").Value = "Retour is nog niet verwerkt"
$ws.Range("B31").Value = "mailmind.test@zohomail.eu"
$ws.Range("C31").Value = "Ik heb iets teruggestuurd maar hoor niks. Wanneer krijg ik mijn geld terug?"
$ws.Range("D31").Value = "Retour / Terugbetaling"
$ws.Range("E31").Value = "Beste klant,`nBedankt voor je bericht. Om je vraag over de terugbetaling te kunnen beantwoorden, hebben we wat meer informatie nodig. Zou je ons alsjeblieft je ordernummer kunnen sturen, zodat we dit kunnen nakijken in ons systeem? Op die manier kunnen we controleren of de retourzending is ontvangen en de terugbetaling is verwerkt.`nWe doen ons best om je zo snel mogelijk van dienst te zijn. Bedankt voor je geduld en medewerking.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("F31").Value = "2025-06-22 19:01:16"
$ws.Range("G31").Value = "Ja"

# Re-measure the row height automatically (keeps it on the default/standard
# height, same as the other rows, instead of leaving a stale custom height
# behind after writing the multi-line "Antwoord" text).
$ws.Rows.Item(31).AutoFit()

# Extend the two conditional-formatting ranges (Categorie / Beantwoord columns)
# so they keep covering the whole table now that it has grown by one row.
$catRules = $ws.Range("D2:D30").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($ws.Range("D2:D31"))
}

$answeredRules = $ws.Range("G2:G30").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($ws.Range("G2:G31"))
}

# --- Dashboard sheet: bump the "Retour / Terugbetaling" tally to match the new log entry ---
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("B4").Value = 4
